$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6130.3335
$ws.Range("I18").Value = 2945.5
$ws.Range("K18").Value = 2945.5
$ws.Range("M18").Value = -2661.5
$ws.Range("H43").Value = 12002.368
$ws.Range("I43").Value = 6999.75
$ws.Range("J43").Value = 13336.4
$ws.Range("K43").Value = 6999.75
$ws.Range("L43").Value = 13336.4
$ws.Range("M43").Value = -6930.75
$ws.Range("N43").Value = -13474.4
$ws.Range("H51").Value = 9500.333000000001
$ws.Range("J51").Value = 9500.333000000001
$ws.Range("L51").Value = 9500.333000000001
$ws.Range("N51").Value = -10468.333
$ws.Range("H61").Value = 195.18182
$ws.Range("I61").Value = 195.18182
$ws.Range("K61").Value = 585.5454599999999
$ws.Range("M61").Value = -413.5454599999999
$ws.Range("H62").Value = 8840.6
$ws.Range("J62").Value = 9715.143
$ws.Range("L62").Value = 9715.143
$ws.Range("N62").Value = -10963.143
$ws.Range("H64").Value = 9614.714
$ws.Range("J64").Value = 10633.833
$ws.Range("L64").Value = 10633.833
$ws.Range("N64").Value = -11129.833
$ws.Range("H65").Value = 8840.6
$ws.Range("J65").Value = 9715.143
$ws.Range("L65").Value = 48575.715
$ws.Range("N65").Value = -54815.715
$ws.Range("H67").Value = 9614.714
$ws.Range("J67").Value = 10633.833
$ws.Range("L67").Value = 10633.833
$ws.Range("N67").Value = -12349.833
$ws.Range("H74").Value = 8641.294
$ws.Range("I74").Value = 6909.8
$ws.Range("J74").Value = 11114.857
$ws.Range("K74").Value = 6909.8
$ws.Range("L74").Value = 11114.857
$ws.Range("M74").Value = -5973.8
$ws.Range("N74").Value = -12986.857
$ws.Range("H77").Value = 8641.294
$ws.Range("I77").Value = 6909.8
$ws.Range("J77").Value = 11114.857
$ws.Range("K77").Value = 34549
$ws.Range("L77").Value = 55574.285
$ws.Range("M77").Value = -29869
$ws.Range("N77").Value = -64934.285
$ws.Range("H92").Value = 4247.1177
$ws.Range("I92").Value = 3330.4814
$ws.Range("J92").Value = 7782.7144
$ws.Range("K92").Value = 3330.4814
$ws.Range("L92").Value = 7782.7144
$ws.Range("M92").Value = -2082.4814
$ws.Range("N92").Value = -10278.7144
$ws.Range("H98").Value = 11093.9
$ws.Range("I98").Value = 1277.4286
$ws.Range("K98").Value = 1277.4286
$ws.Range("M98").Value = 220.5714
$ws.Range("H122").Value = 11093.9
$ws.Range("I122").Value = 1277.4286
$ws.Range("K122").Value = 3832.2858
$ws.Range("M122").Value = -1382.2858
$ws.Range("H123").Value = 66984
$ws.Range("J123").Value = 66984
$ws.Range("L123").Value = 66984
$ws.Range("N123").Value = -76784
$ws.Range("H125").Value = 1813.6
$ws.Range("I125").Value = 1566
$ws.Range("K125").Value = 14094
$ws.Range("M125").Value = -11634
$ws.Range("H132").Value = 2061.3408
$ws.Range("I132").Value = 1642.6842
$ws.Range("J132").Value = 4712.8335
$ws.Range("K132").Value = 4928.0526
$ws.Range("L132").Value = 14138.5005
$ws.Range("M132").Value = -2398.0526
$ws.Range("N132").Value = -19198.5005
$ws.Range("H133").Value = 69997
$ws.Range("J133").Value = 69997
$ws.Range("L133").Value = 69997
$ws.Range("N133").Value = -80117
$ws.Range("H134").Value = 69995.60000000001
$ws.Range("J134").Value = 69995.60000000001
$ws.Range("L134").Value = 69995.60000000001
$ws.Range("N134").Value = -80135.60000000001
$ws.Range("H136").Value = 68993.14
$ws.Range("J136").Value = 68993.14
$ws.Range("L136").Value = 68993.14
$ws.Range("N136").Value = -79193.14
$ws.Range("H137").Value = 3977.0435
$ws.Range("I137").Value = 3239.5557
$ws.Range("J137").Value = 4451.143
$ws.Range("K137").Value = 9718.667099999999
$ws.Range("L137").Value = 13353.429
$ws.Range("M137").Value = -7168.667099999999
$ws.Range("N137").Value = -18453.429
$ws.Range("H139").Value = 65553.22
$ws.Range("J139").Value = 69997.375
$ws.Range("L139").Value = 69997.375
$ws.Range("N139").Value = -80277.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11341
$ws.Range("I2").Value = 2805.8572
$ws.Range("K2").Value = 2805.8572
$ws.Range("M2").Value = -2692.8572
$ws.Range("H32").Value = 2977.9614
$ws.Range("I32").Value = 3541.6843
$ws.Range("J32").Value = 1447.8572
$ws.Range("K32").Value = 3541.6843
$ws.Range("L32").Value = 1447.8572
$ws.Range("M32").Value = -3254.6843
$ws.Range("N32").Value = -2021.8572
$ws.Range("H45").Value = 76926850
$ws.Range("I45").Value = 142857940
$ws.Range("K45").Value = 142857940
$ws.Range("M45").Value = -142857563
$ws.Range("H61").Value = 6073.5186
$ws.Range("I61").Value = 5455.8696
$ws.Range("K61").Value = 5455.8696
$ws.Range("M61").Value = -5243.8696
$ws.Range("H74").Value = 9807199
$ws.Range("I74").Value = 12347817
$ws.Range("J74").Value = 7670
$ws.Range("K74").Value = 12347817
$ws.Range("L74").Value = 7670
$ws.Range("M74").Value = -12346943
$ws.Range("N74").Value = -9418
$ws.Range("H77").Value = 9807199
$ws.Range("I77").Value = 12347817
$ws.Range("J77").Value = 7670
$ws.Range("K77").Value = 61739085
$ws.Range("L77").Value = 38350
$ws.Range("M77").Value = -61734717
$ws.Range("N77").Value = -47086
$ws.Range("H86").Value = 55000
$ws.Range("J86").Value = 60000
$ws.Range("L86").Value = 60000
$ws.Range("N86").Value = -62372
$ws.Range("H89").Value = 55000
$ws.Range("J89").Value = 60000
$ws.Range("L89").Value = 60000
$ws.Range("N89").Value = -191856
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H116").Value = 11341
$ws.Range("I116").Value = 2805.8572
$ws.Range("K116").Value = 2805.8572
$ws.Range("M116").Value = -511.8571999999999
$ws.Range("H132").Value = 6376.4346
$ws.Range("I132").Value = 5299.4326
$ws.Range("K132").Value = 15898.2978
$ws.Range("M132").Value = -13368.2978
$ws.Range("H136").Value = 6073.5186
$ws.Range("I136").Value = 5455.8696
$ws.Range("K136").Value = 16367.6088
$ws.Range("M136").Value = -13817.6088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11341
$ws.Range("I3").Value = 2805.8572
$ws.Range("K3").Value = 2805.8572
$ws.Range("M3").Value = -2691.8572
$ws.Range("H16").Value = 3002.3333
$ws.Range("I16").Value = 3503.5
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 3503.5
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -3333.5
$ws.Range("N16").Value = -2340
$ws.Range("H20").Value = 1781.5714
$ws.Range("I20").Value = 1937.5454
$ws.Range("J20").Value = 1680.6471
$ws.Range("K20").Value = 1937.5454
$ws.Range("L20").Value = 1680.6471
$ws.Range("M20").Value = -1690.5454
$ws.Range("N20").Value = -2174.6471
$ws.Range("H86").Value = 3457.577
$ws.Range("I86").Value = 2332.65
$ws.Range("K86").Value = 2332.65
$ws.Range("M86").Value = -1209.65
$ws.Range("H89").Value = 3457.577
$ws.Range("I89").Value = 2332.65
$ws.Range("K89").Value = 11663.25
$ws.Range("M89").Value = -6047.25
$ws.Range("H94").Value = 1119.7307
$ws.Range("I94").Value = 1029.2916
$ws.Range("J94").Value = 2205
$ws.Range("K94").Value = 1029.2916
$ws.Range("L94").Value = 2205
$ws.Range("M94").Value = -578.2916
$ws.Range("N94").Value = -3107
$ws.Range("H105").Value = 17961.5
$ws.Range("I105").Value = 25209.445
$ws.Range("J105").Value = 8642.714
$ws.Range("K105").Value = 25209.445
$ws.Range("L105").Value = 8642.714
$ws.Range("M105").Value = -23462.445
$ws.Range("N105").Value = -12136.714
$ws.Range("H127").Value = 55000
$ws.Range("J127").Value = 55000
$ws.Range("L127").Value = 55000
$ws.Range("N127").Value = -64920
$ws.Range("H134").Value = 4601.4707
$ws.Range("I134").Value = 2682.5557
$ws.Range("J134").Value = 6760.25
$ws.Range("K134").Value = 8047.6671
$ws.Range("L134").Value = 20280.75
$ws.Range("M134").Value = -5512.6671
$ws.Range("N134").Value = -25350.75
$ws.Range("H140").Value = 66212
$ws.Range("J140").Value = 66212
$ws.Range("L140").Value = 66212
$ws.Range("N140").Value = -76572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2617.9092
$ws.Range("I22").Value = 2084.3845
$ws.Range("J22").Value = 3388.5557
$ws.Range("K22").Value = 2084.3845
$ws.Range("L22").Value = 3388.5557
$ws.Range("M22").Value = -1734.3845
$ws.Range("N22").Value = -4088.5557
$ws.Range("H31").Value = 23164.107
$ws.Range("I31").Value = 2470.9656
$ws.Range("K31").Value = 2470.9656
$ws.Range("M31").Value = -2175.9656
$ws.Range("H34").Value = 23164.107
$ws.Range("I34").Value = 2470.9656
$ws.Range("K34").Value = 2470.9656
$ws.Range("M34").Value = -2268.9656
$ws.Range("H62").Value = 8164.8667
$ws.Range("I62").Value = 5829.3335
$ws.Range("K62").Value = 5829.3335
$ws.Range("M62").Value = -5205.3335
$ws.Range("H65").Value = 8164.8667
$ws.Range("I65").Value = 5829.3335
$ws.Range("K65").Value = 29146.6675
$ws.Range("M65").Value = -26026.6675
$ws.Range("H99").Value = 2626.4119
$ws.Range("I99").Value = 2400
$ws.Range("K99").Value = 2400
$ws.Range("M99").Value = -902
$ws.Range("H105").Value = 4077.6924
$ws.Range("I105").Value = 6349.5
$ws.Range("K105").Value = 6349.5
$ws.Range("M105").Value = -4602.5
$ws.Range("H126").Value = 2626.4119
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730
$ws.Range("H132").Value = 3619.4546
$ws.Range("I132").Value = 2566.1177
$ws.Range("K132").Value = 7698.353099999999
$ws.Range("M132").Value = -5168.353099999999
$ws.Range("H134").Value = 3315.8235
$ws.Range("I134").Value = 1945.6
$ws.Range("K134").Value = 5836.799999999999
$ws.Range("M134").Value = -3301.799999999999
$ws.Range("H135").Value = 69995.8
$ws.Range("J135").Value = 69995.8
$ws.Range("L135").Value = 69995.8
$ws.Range("N135").Value = -80135.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2698585.2
$ws.Range("I4").Value = 2949399.8
$ws.Range("J4").Value = 1512915.9
$ws.Range("K4").Value = 8848199.399999999
$ws.Range("L4").Value = 4538747.699999999
$ws.Range("M4").Value = -8848087.399999999
$ws.Range("N4").Value = -4538971.699999999
$ws.Range("H46").Value = 636.5
$ws.Range("I46").Value = 330.72726
$ws.Range("K46").Value = 992.18178
$ws.Range("M46").Value = -901.18178
$ws.Range("H128").Value = 1583666.6
$ws.Range("I128").Value = 1583666.6
$ws.Range("K128").Value = 4750999.800000001
$ws.Range("M128").Value = -4746019.800000001
$ws.Range("H131").Value = 5721155.5
$ws.Range("I131").Value = 1659.619
$ws.Range("J131").Value = 14960341
$ws.Range("K131").Value = 4978.857
$ws.Range("L131").Value = 44881023
$ws.Range("M131").Value = 61.14300000000003
$ws.Range("N131").Value = -44891103
$ws.Range("H141").Value = 4224.6924
$ws.Range("I141").Value = 447.36365
$ws.Range("K141").Value = 1342.09095
$ws.Range("M141").Value = 3837.90905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.40741
$ws.Range("I2").Value = 60.444443
$ws.Range("J2").Value = 147.33333
$ws.Range("K2").Value = 60.444443
$ws.Range("L2").Value = 147.33333
$ws.Range("M2").Value = 52.555557
$ws.Range("N2").Value = -373.33333
$ws.Range("H15").Value = 50399.6
$ws.Range("J15").Value = 50399.6
$ws.Range("L15").Value = 50399.6
$ws.Range("N15").Value = -50975.6
$ws.Range("H70").Value = 6482.6665
$ws.Range("I70").Value = 6499
$ws.Range("J70").Value = 6450
$ws.Range("K70").Value = 6499
$ws.Range("L70").Value = 6450
$ws.Range("M70").Value = -6229
$ws.Range("N70").Value = -6990
$ws.Range("H73").Value = 6482.6665
$ws.Range("I73").Value = 6499
$ws.Range("J73").Value = 6450
$ws.Range("K73").Value = 6499
$ws.Range("L73").Value = 6450
$ws.Range("M73").Value = -5563
$ws.Range("N73").Value = -8322
$ws.Range("H81").Value = 50399.6
$ws.Range("J81").Value = 50399.6
$ws.Range("L81").Value = 50399.6
$ws.Range("N81").Value = -52395.6
$ws.Range("H84").Value = 50399.6
$ws.Range("J84").Value = 50399.6
$ws.Range("L84").Value = 151198.8
$ws.Range("N84").Value = -161182.8
$ws.Range("H97").Value = 1975.9333
$ws.Range("I97").Value = 1440
$ws.Range("J97").Value = 3449.75
$ws.Range("K97").Value = 1440
$ws.Range("L97").Value = 3449.75
$ws.Range("M97").Value = -944
$ws.Range("N97").Value = -4441.75
$ws.Range("H122").Value = 10477
$ws.Range("I122").Value = 6499
$ws.Range("J122").Value = 10725.625
$ws.Range("K122").Value = 19497
$ws.Range("L122").Value = 32176.875
$ws.Range("M122").Value = -17047
$ws.Range("N122").Value = -37076.875
$ws.Range("H133").Value = 69997
$ws.Range("J133").Value = 69997
$ws.Range("L133").Value = 69997
$ws.Range("N133").Value = -80117
$ws.Range("H137").Value = 74999
$ws.Range("J137").Value = 74999
$ws.Range("L137").Value = 74999
$ws.Range("N137").Value = -85199
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 74998
$ws.Range("J140").Value = 74998
$ws.Range("L140").Value = 74998
$ws.Range("N140").Value = -85358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3720.8572
$ws.Range("I46").Value = 915.3333
$ws.Range("J46").Value = 4843.067
$ws.Range("K46").Value = 915.3333
$ws.Range("L46").Value = 4843.067
$ws.Range("M46").Value = -727.3333
$ws.Range("N46").Value = -5219.067
$ws.Range("H60").Value = 25500
$ws.Range("J60").Value = 25500
$ws.Range("L60").Value = 25500
$ws.Range("N60").Value = -26518
$ws.Range("H82").Value = 3865.7693
$ws.Range("I82").Value = 3426.9375
$ws.Range("J82").Value = 4567.9
$ws.Range("K82").Value = 3426.9375
$ws.Range("L82").Value = 4567.9
$ws.Range("M82").Value = -3065.9375
$ws.Range("N82").Value = -5289.9
$ws.Range("H85").Value = 3865.7693
$ws.Range("I85").Value = 3426.9375
$ws.Range("J85").Value = 4567.9
$ws.Range("K85").Value = 3426.9375
$ws.Range("L85").Value = 4567.9
$ws.Range("M85").Value = -2178.9375
$ws.Range("N85").Value = -7063.9
$ws.Range("H97").Value = 24500
$ws.Range("J97").Value = 24500
$ws.Range("L97").Value = 24500
$ws.Range("N97").Value = -26482
$ws.Range("H132").Value = 4913.037
$ws.Range("I132").Value = 3877.9443
$ws.Range("K132").Value = 11633.8329
$ws.Range("M132").Value = -9103.832900000001
$ws.Range("H133").Value = 97659.836
$ws.Range("J133").Value = 97659.836
$ws.Range("L133").Value = 97659.836
$ws.Range("N133").Value = -102719.836
$ws.Range("H134").Value = 95000
$ws.Range("J134").Value = 95000
$ws.Range("L134").Value = 95000
$ws.Range("N134").Value = -105140
$ws.Range("H136").Value = 6671.85
$ws.Range("I136").Value = 4782.5454
$ws.Range("J136").Value = 8981
$ws.Range("K136").Value = 14347.6362
$ws.Range("L136").Value = 26943
$ws.Range("M136").Value = -11797.6362
$ws.Range("N136").Value = -32043
$ws.Range("H137").Value = 70903.5
$ws.Range("J137").Value = 70903.5
$ws.Range("L137").Value = 70903.5
$ws.Range("N137").Value = -81103.5
$ws.Range("H139").Value = 67371.25
$ws.Range("J139").Value = 67371.25
$ws.Range("L139").Value = 67371.25
$ws.Range("N139").Value = -77651.25
$ws.Range("H141").Value = 78299.2
$ws.Range("J141").Value = 78299.2
$ws.Range("L141").Value = 78299.2
$ws.Range("N141").Value = -88659.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4105.2104
$ws.Range("I126").Value = 2900.3333
$ws.Range("K126").Value = 8700.999899999999
$ws.Range("M126").Value = -6230.999899999999
$ws.Range("H132").Value = 3589.2856
$ws.Range("I132").Value = 2402.2793
$ws.Range("K132").Value = 7206.8379
$ws.Range("M132").Value = -4676.8379
$ws.Range("H135").Value = 69995
$ws.Range("J135").Value = 69995
$ws.Range("L135").Value = 69995
$ws.Range("N135").Value = -80135
$ws.Range("H136").Value = 3555.76
$ws.Range("I136").Value = 2022.9412
$ws.Range("K136").Value = 6068.8236
$ws.Range("M136").Value = -3518.8236
$ws.Range("H137").Value = 65423
$ws.Range("J137").Value = 65423
$ws.Range("L137").Value = 65423
$ws.Range("N137").Value = -75623
$ws.Range("H139").Value = 69997
$ws.Range("J139").Value = 69997
$ws.Range("L139").Value = 69997
$ws.Range("N139").Value = -80277
$ws.Range("H141").Value = 154140.6
$ws.Range("J141").Value = 154140.6
$ws.Range("L141").Value = 154140.6
$ws.Range("N141").Value = -164500.6
